$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.36
$ws.Range("C2").Value = 43.58
$ws.Range("D2").Value = 8.539999999999999
$ws.Range("E2").Value = 84.97
$ws.Range("F2").Value = 55.89
$ws.Range("G2").Value = 95.84999999999999
$ws.Range("H2").Value = 89.02
$ws.Range("I2").Value = 55.56
$ws.Range("J2").Value = 46.47
$ws.Range("K2").Value = 19.76
$ws.Range("L2").Value = 9.26
$ws.Range("M2").Value = 63.8
$ws.Range("N2").Value = 16.72
$ws.Range("O2").Value = 56.57
$ws.Range("P2").Value = 98.17

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 27.33
$ws.Range("D3").Value = 2.63
$ws.Range("E3").Value = 90.86
$ws.Range("F3").Value = 49.42
$ws.Range("G3").Value = 98.83
$ws.Range("H3").Value = 93.68000000000001
$ws.Range("I3").Value = 53.79
$ws.Range("J3").Value = 27.14
$ws.Range("K3").Value = 25
$ws.Range("L3").Value = 4.47
$ws.Range("M3").Value = 36.38
$ws.Range("N3").Value = 12.06
$ws.Range("O3").Value = 53.6
$ws.Range("P3").Value = 93.39

# Row 4
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 5.88
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 94.12
$ws.Range("F4").Value = 85.29000000000001
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 97.06
$ws.Range("I4").Value = 20.59
$ws.Range("J4").Value = 67.65000000000001
$ws.Range("K4").Value = 52.94
$ws.Range("L4").Value = 32.35
$ws.Range("M4").Value = 82.34999999999999
$ws.Range("N4").Value = 58.82
$ws.Range("O4").Value = 20.59
$ws.Range("P4").Value = 88.23999999999999

# Row 5 - clear the values (becomes empty inline strings)
$ws.Range("B5:P5").ClearContents()
